# Weekly fruit/vegetable price update:
# Two new weekly records (date 2023-06-16, serial 45093) are added at the
# top of the "Caqui" (Mankaki variety) price series, pushing the existing
# rows 65-79 down to rows 67-81.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows right above the current first data row of the block (row 65).
$ws.Rows("65:66").Insert()

# --- New row 65: "Primera" quality ---
$ws.Range("A65").Value = 5
$ws.Range("B65").Value = "Macroferia Regional de Talca"
$ws.Range("C65").Value = "Maule"
$ws.Range("D65").Value = 45093
$ws.Range("E65").Value = 7
$ws.Range("F65").Value = "Fruta"
$ws.Range("G65").Value = 100107
$ws.Range("H65").Value = "Otros"
$ws.Range("I65").Value = 100107001
$ws.Range("J65").Value = "Caqui"
$ws.Range("K65").Value = "Mankaki"
$ws.Range("L65").Value = "Primera"
$ws.Range("M65").Value = 50
$ws.Range("N65").Value = 18000
$ws.Range("O65").Value = 18000
$ws.Range("P65").Value = 18000
$ws.Range("Q65").Value = "$/caja 18 kilos granel"
$ws.Range("R65").Value = "Región del Maule"
$ws.Range("S65").Value = 1000
$ws.Range("T65").Value = 18

# --- New row 66: "Segunda" quality ---
$ws.Range("A66").Value = 5
$ws.Range("B66").Value = "Macroferia Regional de Talca"
$ws.Range("C66").Value = "Maule"
$ws.Range("D66").Value = 45093
$ws.Range("E66").Value = 7
$ws.Range("F66").Value = "Fruta"
$ws.Range("G66").Value = 100107
$ws.Range("H66").Value = "Otros"
$ws.Range("I66").Value = 100107001
$ws.Range("J66").Value = "Caqui"
$ws.Range("K66").Value = "Mankaki"
$ws.Range("L66").Value = "Segunda"
$ws.Range("M66").Value = 40
$ws.Range("N66").Value = 15000
$ws.Range("O66").Value = 15000
$ws.Range("P66").Value = 15000
$ws.Range("Q66").Value = "$/caja 18 kilos granel"
$ws.Range("R66").Value = "Región del Maule"
$ws.Range("S66").Value = 833
$ws.Range("T66").Value = 18
